$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 2-28) currently stores the exercicio/year as a plain
# number (2015, 2016, ... 2023, repeated for each account block). The
# edit replaces those numeric years with literal text date strings
# "01/01/<year>" (stored as shared strings, not as Excel date serials).
$years = @(2015, 2016, 2017, 2018, 2019, 2020, 2021, 2022, 2023, `
           2015, 2016, 2017, 2018, 2019, 2020, 2021, 2022, 2023, `
           2015, 2016, 2017, 2018, 2019, 2020, 2021, 2022, 2023)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)

    # Force text storage so Excel doesn't reinterpret the "dd/mm/yyyy"
    # looking string as a date serial number.
    $cell.NumberFormat = "@"
    $cell.Value = "01/01/" + $years[$i]
    # Drop back to the default (unstyled) cell format - the source
    # cells had no explicit style and the edit doesn't introduce one.
    $cell.Style = "Normal"
}
